$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.361.56"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "3.268.54"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "615.65"
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("D6").Value = "157.80"
$ws.Range("E6").Value = "  +1.93%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.265.68"
$ws.Range("E8").Value = "  +2.72%  "
$ws.Range("D9").Value = "0.545"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("E10").Value = "  +2.00%  "
$ws.Range("D11").Value = "5.79"
$ws.Range("E11").Value = "  +1.64%  "
$ws.Range("D12").Value = "0.495"
$ws.Range("E12").Value = "  -4.21%  "
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("D14").Value = "39.06"
$ws.Range("E14").Value = "  +2.03%  "
$ws.Range("D15").Value = "3.805.52"
$ws.Range("E15").Value = "  +3.01%  "
$ws.Range("D16").Value = "66.453.90"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").Value = "7.44"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "3.262.15"
$ws.Range("E18").Value = "  +2.61%  "
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("D20").Value = "505.50"
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").Value = "15.52"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("D22").Value = "0.756"
$ws.Range("E22").Value = "  +3.70%  "
$ws.Range("D24").Value = "14.65"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("D25").Value = "87.05"
$ws.Range("E25").Value = "  +3.09%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "3.06"
$ws.Range("E27").Value = "  +1.92%  "
$ws.Range("D28").Value = "9.24"
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("E30").Value = "  +45.33%  "
$ws.Range("D31").Value = "7.01"
$ws.Range("E31").Value = "  -2.40%  "
$ws.Range("E32").Value = "  -3.94%  "
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  -2.92%  "
$ws.Range("D36").Value = "6.47"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D37").Value = "3.43"
$ws.Range("E37").Value = "  +20.20%  "
$ws.Range("D38").Value = "55.62"
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("D39").Value = "0.0₃0792"
$ws.Range("E39").Value = "  +15.59%  "
$ws.Range("D40").Value = "494.90"
$ws.Range("E40").Value = "  -2.03%  "
$ws.Range("D41").Value = "0.0425"
$ws.Range("E41").Value = "  +1.46%  "
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("D43").Value = "8.83"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("D44").Value = "2.54"
$ws.Range("E44").Value = "  +3.81%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "3.014.02"
$ws.Range("E45").Value = "  +6.52%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "0.293"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("D47").Value = "28.97"
$ws.Range("E47").Value = "  +3.33%  "
$ws.Range("E48").Value = "  +6.35%  "
$ws.Range("E49").Value = "  +2.51%  "
$ws.Range("D51").Value = "2.52"
$ws.Range("E51").Value = "  -3.47%  "
